# Update cryptos list price (D) and 1h volume % (E) columns per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.163.52"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "'2.933.45"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'593.40"
$ws.Range("D6").Value = "'145.17"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'7.00"
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "'33.78"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "'3.419.80"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "'61.100.90"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "'6.74"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'2.932.23"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "'433.88"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'13.51"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'7.12"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'81.72"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'11.09"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "'11.89"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  +3.20%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "'1.02"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").Value = "'2.98"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'8.62"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'42.08"
$ws.Range("E41").Value = "  +5.05%  "
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").Value = "'374.91"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "'2.708.50"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "'133.81"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'23.94"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'0.125"
$ws.Range("E51").Value = "  +0.13%  "
